# Línea 141 - refresh de horarios (scrap 12:54:24).
# Actualiza las 3 hojas (LP1912, LP1912-215, 6203-6173): re-escribe filas
# existentes con los nuevos valores scrapeados y agrega las filas nuevas
# al final de cada tabla.
$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:54:24"
$ws.Cells.Item(3, 1).Value = "Total filas: 237"
$ws.Cells.Item(47, 3).Value = "14_ABASTO"
$ws.Cells.Item(48, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(63, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(64, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(75, 3).Value = "10_OLMOS"
$ws.Cells.Item(76, 3).Value = "215D_EL PATO"
$ws.Cells.Item(88, 1).Value = "08:13:38"
$ws.Cells.Item(88, 3).Value = "215B_EL PATO"
$ws.Cells.Item(88, 4).Value = 10
$ws.Cells.Item(89, 1).Value = "07:28:14"
$ws.Cells.Item(89, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(89, 4).Value = 55
$ws.Cells.Item(99, 1).Value = "08:13:38"
$ws.Cells.Item(99, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(99, 4).Value = 48
$ws.Cells.Item(100, 1).Value = "08:56:26"
$ws.Cells.Item(100, 3).Value = "215A_EL PATO"
$ws.Cells.Item(100, 4).Value = 5
$ws.Cells.Item(113, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(114, 1).Value = "08:56:26"
$ws.Cells.Item(114, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(114, 4).Value = 27
$ws.Cells.Item(115, 1).Value = "07:28:14"
$ws.Cells.Item(115, 3).Value = "17_ROMERO"
$ws.Cells.Item(115, 4).Value = 115
$ws.Cells.Item(120, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(121, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(163, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(164, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(198, 3).Value = "17_179 Y 38"
$ws.Cells.Item(199, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(207, 1).Value = "12:54:24"
$ws.Cells.Item(207, 2).Value = "12:54"
$ws.Cells.Item(207, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(207, 4).Value = 0
$ws.Cells.Item(208, 1).Value = "12:54:24"
$ws.Cells.Item(208, 2).Value = "12:55"
$ws.Cells.Item(208, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(208, 4).Value = 1
$ws.Cells.Item(209, 2).Value = "12:55"
$ws.Cells.Item(209, 3).Value = "10_OLMOS"
$ws.Cells.Item(209, 4).Value = 18
$ws.Cells.Item(210, 2).Value = "13:02"
$ws.Cells.Item(210, 3).Value = "15_ABASTO"
$ws.Cells.Item(210, 4).Value = 25
$ws.Cells.Item(211, 1).Value = "12:37:00"
$ws.Cells.Item(211, 2).Value = "13:03"
$ws.Cells.Item(211, 3).Value = "14_ABASTO"
$ws.Cells.Item(211, 4).Value = 26
$ws.Cells.Item(212, 1).Value = "12:54:24"
$ws.Cells.Item(212, 2).Value = "13:05"
$ws.Cells.Item(212, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(212, 4).Value = 11
$ws.Cells.Item(213, 2).Value = "13:06"
$ws.Cells.Item(213, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(213, 4).Value = 29
$ws.Cells.Item(214, 1).Value = "12:54:24"
$ws.Cells.Item(214, 2).Value = "13:06"
$ws.Cells.Item(214, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(215, 1).Value = "12:37:00"
$ws.Cells.Item(215, 2).Value = "13:07"
$ws.Cells.Item(215, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(215, 4).Value = 30
$ws.Cells.Item(216, 1).Value = "12:54:24"
$ws.Cells.Item(216, 2).Value = "13:08"
$ws.Cells.Item(216, 3).Value = "10_OLMOS"
$ws.Cells.Item(216, 4).Value = 14
$ws.Cells.Item(217, 1).Value = "11:22:44"
$ws.Cells.Item(217, 2).Value = "13:10"
$ws.Cells.Item(217, 3).Value = "10_OLMOS"
$ws.Cells.Item(217, 4).Value = 108
$ws.Cells.Item(218, 1).Value = "11:22:44"
$ws.Cells.Item(218, 2).Value = "13:13"
$ws.Cells.Item(218, 3).Value = "215D_EL PATO"
$ws.Cells.Item(218, 4).Value = 111
$ws.Cells.Item(219, 1).Value = "12:54:24"
$ws.Cells.Item(219, 2).Value = "13:14"
$ws.Cells.Item(219, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(219, 4).Value = 20
$ws.Cells.Item(220, 1).Value = "12:54:24"
$ws.Cells.Item(220, 2).Value = "13:14"
$ws.Cells.Item(220, 3).Value = "215D_EL PATO"
$ws.Cells.Item(220, 4).Value = 20
$ws.Cells.Item(221, 1).Value = "12:37:00"
$ws.Cells.Item(221, 2).Value = "13:15"
$ws.Cells.Item(221, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(221, 4).Value = 38
$ws.Cells.Item(222, 1).Value = "12:54:24"
$ws.Cells.Item(222, 2).Value = "13:19"
$ws.Cells.Item(222, 3).Value = "10_OLMOS"
$ws.Cells.Item(222, 4).Value = 25
$ws.Cells.Item(223, 2).Value = "13:20"
$ws.Cells.Item(223, 3).Value = "10_OLMOS"
$ws.Cells.Item(223, 4).Value = 43
$ws.Cells.Item(224, 1).Value = "12:54:24"
$ws.Cells.Item(224, 2).Value = "13:21"
$ws.Cells.Item(224, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(224, 4).Value = 27
$ws.Cells.Item(225, 2).Value = "13:23"
$ws.Cells.Item(225, 3).Value = "10_OLMOS"
$ws.Cells.Item(225, 4).Value = 81
$ws.Cells.Item(226, 1).Value = "12:54:24"
$ws.Cells.Item(226, 2).Value = "13:26"
$ws.Cells.Item(226, 3).Value = "15_ABASTO"
$ws.Cells.Item(226, 4).Value = 32
$ws.Cells.Item(227, 1).Value = "12:54:24"
$ws.Cells.Item(227, 2).Value = "13:26"
$ws.Cells.Item(227, 3).Value = "14_ABASTO"
$ws.Cells.Item(227, 4).Value = 32
$ws.Cells.Item(228, 1).Value = "12:37:00"
$ws.Cells.Item(228, 2).Value = "13:27"
$ws.Cells.Item(228, 3).Value = "14_ABASTO"
$ws.Cells.Item(228, 4).Value = 50
$ws.Cells.Item(229, 1).Value = "12:54:24"
$ws.Cells.Item(229, 2).Value = "13:34"
$ws.Cells.Item(229, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(229, 4).Value = 40
$ws.Cells.Item(230, 1).Value = "12:54:24"
$ws.Cells.Item(230, 2).Value = "13:36"
$ws.Cells.Item(230, 3).Value = "15_ABASTO"
$ws.Cells.Item(230, 4).Value = 42
$ws.Cells.Item(231, 1).Value = "12:54:24"
$ws.Cells.Item(231, 2).Value = "13:46"
$ws.Cells.Item(231, 3).Value = "17_ROMERO"
$ws.Cells.Item(231, 4).Value = 52
$ws.Cells.Item(232, 1).Value = "12:54:24"
$ws.Cells.Item(232, 2).Value = "13:50"
$ws.Cells.Item(232, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(232, 4).Value = 56
$ws.Cells.Item(233, 1).Value = "12:54:24"
$ws.Cells.Item(233, 2).Value = "13:50"
$ws.Cells.Item(233, 3).Value = "215A_EL PATO"
$ws.Cells.Item(233, 4).Value = 56
$ws.Cells.Item(234, 1).Value = "12:37:00"
$ws.Cells.Item(234, 2).Value = "13:51"
$ws.Cells.Item(234, 3).Value = "215A_EL PATO"
$ws.Cells.Item(234, 4).Value = 74
$ws.Cells.Item(234, 5).Value = "LP1912"
$ws.Cells.Item(235, 1).Value = "12:54:24"
$ws.Cells.Item(235, 2).Value = "13:56"
$ws.Cells.Item(235, 3).Value = "225_GOMEZ"
$ws.Cells.Item(235, 4).Value = 62
$ws.Cells.Item(235, 5).Value = "LP1912"
$ws.Cells.Item(236, 1).Value = "12:54:24"
$ws.Cells.Item(236, 2).Value = "13:56"
$ws.Cells.Item(236, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(236, 4).Value = 62
$ws.Cells.Item(236, 5).Value = "LP1912"
$ws.Cells.Item(237, 1).Value = "12:37:00"
$ws.Cells.Item(237, 2).Value = "13:57"
$ws.Cells.Item(237, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(237, 4).Value = 80
$ws.Cells.Item(237, 5).Value = "LP1912"
$ws.Cells.Item(238, 1).Value = "12:54:24"
$ws.Cells.Item(238, 2).Value = "14:04"
$ws.Cells.Item(238, 3).Value = "17_ROMERO"
$ws.Cells.Item(238, 4).Value = 70
$ws.Cells.Item(238, 5).Value = "LP1912"
$ws.Cells.Item(239, 1).Value = "12:54:24"
$ws.Cells.Item(239, 2).Value = "14:17"
$ws.Cells.Item(239, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(239, 4).Value = 83
$ws.Cells.Item(239, 5).Value = "LP1912"
$ws.Cells.Item(240, 1).Value = "12:54:24"
$ws.Cells.Item(240, 2).Value = "14:20"
$ws.Cells.Item(240, 3).Value = "215C_EL PATO"
$ws.Cells.Item(240, 4).Value = 86
$ws.Cells.Item(240, 5).Value = "LP1912"
$ws.Cells.Item(241, 1).Value = "12:54:24"
$ws.Cells.Item(241, 2).Value = "14:21"
$ws.Cells.Item(241, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(241, 4).Value = 87
$ws.Cells.Item(241, 5).Value = "LP1912"
$ws.Cells.Item(242, 1).Value = "12:54:24"
$ws.Cells.Item(242, 2).Value = "14:39"
$ws.Cells.Item(242, 3).Value = "14_ABASTO"
$ws.Cells.Item(242, 4).Value = 105
$ws.Cells.Item(242, 5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:54:24"
$ws.Cells.Item(32, 1).Value = "12:54:24"
$ws.Cells.Item(32, 4).Value = 20
$ws.Cells.Item(33, 1).Value = "12:54:24"
$ws.Cells.Item(33, 4).Value = 56
$ws.Cells.Item(35, 1).Value = "12:54:24"
$ws.Cells.Item(35, 4).Value = 86

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:54:24"
$ws.Cells.Item(3, 1).Value = "Total filas: 35"
$ws.Cells.Item(37, 1).Value = "12:54:24"
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(38, 1).Value = "12:54:24"
$ws.Cells.Item(38, 4).Value = 37
$ws.Cells.Item(39, 1).Value = "12:54:24"
$ws.Cells.Item(39, 4).Value = 75
$ws.Cells.Item(40, 1).Value = "12:54:24"
$ws.Cells.Item(40, 2).Value = "14:53"
$ws.Cells.Item(40, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(40, 4).Value = 119
$ws.Cells.Item(40, 5).Value = "L6203"

